$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fix boot type set, get": the Books sheet's "type" column (D) held
# all-caps genre labels ("FANTASY" / "HORROR"); normalize to title case.
$ws.Range("D2").Value = "Fantasy"
$ws.Range("D3").Value = "Horror"

# Selection ends up on D3 (the last-edited "type" cell).
$ws.Range("D3").Select()
